$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conf")

# Remember the two existing hyperlink target URLs (currently anchored at B2 and B8)
$url1 = "http://bart.ideam.gov.co/cneideam/Capasgeo/"
$url2 = "https://geoportal.dane.gov.co/descargas/veredas/CRVeredas_2017.zip"

# Drop the existing hyperlinks up front; they'll be recreated at their new
# (shifted) locations once the row insert below has moved everything down.
$ws.Hyperlinks.Delete()

# Insert a new row above row 2, shifting existing rows 2-9 down to rows 3-10
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new parameter
$ws.Range("A2").Value = "glo_crs"
$ws.Range("B2").Value = 3116

# Re-create the hyperlinks at their shifted locations: B3 (was B2) and B9 (was B8)
$ws.Hyperlinks.Add($ws.Range("B3"), $url1)
$ws.Range("B3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B9"), $url2)
$ws.Range("B9").Style = "Hyperlink"

# Update the active cell selection to match the post-edit state
$ws.Range("A4").Select()
